$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update mean_fit_time, std_fit_time, mean_score_time, std_score_time
# for rows 2 and 3 with the newly measured values.
$ws.Range("A2").Value = 5.962288459142049
$ws.Range("B2").Value = 0.2950884536539055
$ws.Range("C2").Value = 0.005705515543619792
$ws.Range("D2").Value = 0.0004467918512696173

$ws.Range("A3").Value = 5.663124958674113
$ws.Range("B3").Value = 0.1786572799174025
$ws.Range("C3").Value = 0.005200386047363281
$ws.Range("D3").Value = 0.0002822886716117706
